# Fruta / hortaliza, semanal
# Insert a new week's worth of data (3 rows: Especial / Primera / Segunda for
# Murcott mandarins at Femacal de La Calera) at the top of the data block
# (row 501), pushing the existing rows 501:571 down to 504:574.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before the current row 501 (shifts 501:571 -> 504:574)
$ws.Range("A501:T503").EntireRow.Insert()

# New date for this week's entries: 2021-11-22 (serial 44522)
$newDate = 44522

# Row 501: Murcott - Especial
$ws.Cells.Item(501, 1).Value = 3
$ws.Cells.Item(501, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(501, 3).Value = "Coquimbo"
$ws.Cells.Item(501, 4).Value = $newDate
$ws.Cells.Item(501, 5).Value = 5
$ws.Cells.Item(501, 6).Value = "Fruta"
$ws.Cells.Item(501, 7).Value = 100102
$ws.Cells.Item(501, 8).Value = "Cítricos"
$ws.Cells.Item(501, 9).Value = 100102004
$ws.Cells.Item(501, 10).Value = "Mandarina"
$ws.Cells.Item(501, 11).Value = "Murcott"
$ws.Cells.Item(501, 12).Value = "Especial"
$ws.Cells.Item(501, 13).Value = 65
$ws.Cells.Item(501, 14).Value = 6000
$ws.Cells.Item(501, 15).Value = 6000
$ws.Cells.Item(501, 16).Value = 6000
$ws.Cells.Item(501, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(501, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(501, 19).Value = 600
$ws.Cells.Item(501, 20).Value = 10

# Row 502: Murcott - Primera
$ws.Cells.Item(502, 1).Value = 3
$ws.Cells.Item(502, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(502, 3).Value = "Coquimbo"
$ws.Cells.Item(502, 4).Value = $newDate
$ws.Cells.Item(502, 5).Value = 5
$ws.Cells.Item(502, 6).Value = "Fruta"
$ws.Cells.Item(502, 7).Value = 100102
$ws.Cells.Item(502, 8).Value = "Cítricos"
$ws.Cells.Item(502, 9).Value = 100102004
$ws.Cells.Item(502, 10).Value = "Mandarina"
$ws.Cells.Item(502, 11).Value = "Murcott"
$ws.Cells.Item(502, 12).Value = "Primera"
$ws.Cells.Item(502, 13).Value = 60
$ws.Cells.Item(502, 14).Value = 5000
$ws.Cells.Item(502, 15).Value = 5000
$ws.Cells.Item(502, 16).Value = 5000
$ws.Cells.Item(502, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(502, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(502, 19).Value = 500
$ws.Cells.Item(502, 20).Value = 10

# Row 503: Murcott - Segunda
$ws.Cells.Item(503, 1).Value = 3
$ws.Cells.Item(503, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(503, 3).Value = "Coquimbo"
$ws.Cells.Item(503, 4).Value = $newDate
$ws.Cells.Item(503, 5).Value = 5
$ws.Cells.Item(503, 6).Value = "Fruta"
$ws.Cells.Item(503, 7).Value = 100102
$ws.Cells.Item(503, 8).Value = "Cítricos"
$ws.Cells.Item(503, 9).Value = 100102004
$ws.Cells.Item(503, 10).Value = "Mandarina"
$ws.Cells.Item(503, 11).Value = "Murcott"
$ws.Cells.Item(503, 12).Value = "Segunda"
$ws.Cells.Item(503, 13).Value = 70
$ws.Cells.Item(503, 14).Value = 4000
$ws.Cells.Item(503, 15).Value = 4000
$ws.Cells.Item(503, 16).Value = 4000
$ws.Cells.Item(503, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(503, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(503, 19).Value = 400
$ws.Cells.Item(503, 20).Value = 10

Write-Host "Final dimension:" $ws.UsedRange.Address()
